$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a plain-text value to a cell while preserving its
# original "Text" cell-type (these sheets store numbers/percentages as
# literal strings, not numeric cells) and leaving the cell style unchanged.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '305.79'
Set-TextValue "E2" '-3.79%'
Set-TextValue "D3" '37.05'
Set-TextValue "E3" '-6.94%'
Set-TextValue "D4" '5.101'
Set-TextValue "E4" '-1.03%'
Set-TextValue "D5" '0.07721'
Set-TextValue "E5" '-6.20%'
Set-TextValue "D6" '4.371'
Set-TextValue "E6" '0.77%'
Set-TextValue "D7" '8.198'
Set-TextValue "E7" '-1.85%'
Set-TextValue "D8" '1.881'
Set-TextValue "E8" '-8.69%'
Set-TextValue "D9" '3.194'
Set-TextValue "E9" '-4.27%'
Set-TextValue "D10" '0.9187'
Set-TextValue "E10" '-2.06%'
Set-TextValue "D11" '0.1225'
Set-TextValue "E11" '-10.42%'
Set-TextValue "D12" '0.1898'
Set-TextValue "E12" '-4.21%'
Set-TextValue "D13" '0.08735'
Set-TextValue "E13" '-3.91%'
Set-TextValue "D14" '0.03375'
Set-TextValue "E14" '-3.84%'
Set-TextValue "D15" '0.09702'
Set-TextValue "E15" '-1.13%'
Set-TextValue "D16" '0.001370'
Set-TextValue "E16" '-2.62%'
Set-TextValue "D17" '0.006070'
Set-TextValue "E17" '0.46%'
Set-TextValue "D18" '3.557'
Set-TextValue "E18" '-3.65%'
Set-TextValue "E19" '-3.00%'
Set-TextValue "D20" '0.1283'
Set-TextValue "E20" '-2.08%'
Set-TextValue "D21" '5.029'
Set-TextValue "E21" '1.35%'
Set-TextValue "D22" '0.2499'
Set-TextValue "E22" '2.01%'
Set-TextValue "D23" '0.02113'
Set-TextValue "E23" '5,184.21%'
Set-TextValue "D24" '0.04327'
Set-TextValue "E24" '-0.65%'
Set-TextValue "D25" '0.001218'
Set-TextValue "E25" '-0.84%'
Set-TextValue "D26" '0.004468'
Set-TextValue "E26" '-7.47%'
Set-TextValue "E27" '4.36%'
Set-TextValue "D39" '0.02210'
Set-TextValue "E39" '-1.06%'
Set-TextValue "D40" '0.04922'
Set-TextValue "E40" '-5.44%'
Set-TextValue "E41" '-1.79%'
Set-TextValue "D42" '0.009849'
Set-TextValue "E42" '1.75%'
Set-TextValue "D43" '0.1330'
Set-TextValue "E43" '-5.64%'
Set-TextValue "D44" '0.002002'
Set-TextValue "E44" '-2.22%'
Set-TextValue "D45" '0.008811'
Set-TextValue "E45" '-8.79%'
Set-TextValue "D46" '0.00006821'
Set-TextValue "E46" '2.95%'
Set-TextValue "D47" '0.00000000753'
Set-TextValue "E47" '0.49%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue "D48" '0.003009'
Set-TextValue "E48" '2.28%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue "D49" '0.001306'
Set-TextValue "E49" '-22.69%'
Set-TextValue "D50" '0.00002110'
Set-TextValue "E50" '0.49%'
Set-TextValue "D51" '0.0002009'
Set-TextValue "E51" '0.49%'
